# Update odds and correct-score values in row 4 of Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G4").Value = 2.75
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 2.75
$ws.Range("J4").Value = 3.4
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 2.25
$ws.Range("R4").Value = 1.62
$ws.Range("S4").Value = 4
$ws.Range("T4").Value = 1.22
$ws.Range("U4").Value = 1.5
$ws.Range("V4").Value = 2.5
$ws.Range("W4").Value = 1.83
$ws.Range("X4").Value = 1.83
$ws.Range("Y4").Value = 8
$ws.Range("AC4").Value = 23
$ws.Range("AE4").Value = 7.5
$ws.Range("AF4").Value = 5.5
$ws.Range("AG4").Value = 15
$ws.Range("AH4").Value = 51
$ws.Range("AI4").Value = 301
$ws.Range("AJ4").Value = 8
